$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Corr/total marks
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 45
$ws.Range("E12").Value = "45/140"
